$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1359.8916
$ws.Range("I15").Value = 1359.8916
$ws.Range("K15").Value = 4079.6748
$ws.Range("M15").Value = -3910.6748
$ws.Range("H18").Value = 2683.1667
$ws.Range("I18").Value = 649.5
$ws.Range("K18").Value = 649.5
$ws.Range("M18").Value = -365.5
$ws.Range("H28").Value = 1659.8334
$ws.Range("I28").Value = 1620
$ws.Range("K28").Value = 1620
$ws.Range("M28").Value = -1135
$ws.Range("H41").Value = 1081
$ws.Range("I41").Value = 184.83333
$ws.Range("K41").Value = 184.83333
$ws.Range("M41").Value = 255.16667
$ws.Range("H51").Value = 5099.4443
$ws.Range("H86").Value = 5035.143
$ws.Range("I86").Value = 4707.6665
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 4707.6665
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -3584.6665
$ws.Range("N86").Value = -9246
$ws.Range("H89").Value = 5035.143
$ws.Range("I89").Value = 4707.6665
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 23538.3325
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -17922.3325
$ws.Range("N89").Value = -46232
$ws.Range("H137").Value = 11077646
$ws.Range("I137").Value = 501574.25
$ws.Range("J137").Value = 30306868
$ws.Range("K137").Value = 1504722.75
$ws.Range("L137").Value = 90920604
$ws.Range("M137").Value = -1502172.75
$ws.Range("N137").Value = -90925704
$ws.Range("H138").Value = 3012.11
$ws.Range("I138").Value = 1369.0526
$ws.Range("J138").Value = 3397.5186
$ws.Range("K138").Value = 4107.1578
$ws.Range("L138").Value = 10192.5558
$ws.Range("M138").Value = 1032.8422
$ws.Range("N138").Value = -20472.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12675.791
$ws.Range("I2").Value = 16651.25
$ws.Range("K2").Value = 16651.25
$ws.Range("M2").Value = -16538.25
$ws.Range("H5").Value = 1300
$ws.Range("I5").Value = 1287.5
$ws.Range("K5").Value = 1287.5
$ws.Range("M5").Value = -1175.5
$ws.Range("H32").Value = 15887.068
$ws.Range("I32").Value = 14608.516
$ws.Range("J32").Value = 24979
$ws.Range("K32").Value = 14608.516
$ws.Range("L32").Value = 24979
$ws.Range("M32").Value = -14321.516
$ws.Range("N32").Value = -25553
$ws.Range("H42").Value = 9500
$ws.Range("J42").Value = 9500
$ws.Range("L42").Value = 9500
$ws.Range("N42").Value = -10472
$ws.Range("H53").Value = 4150
$ws.Range("I53").Value = 4150
$ws.Range("K53").Value = 4150
$ws.Range("M53").Value = -3468
$ws.Range("H61").Value = 6238.514
$ws.Range("I61").Value = 5915
$ws.Range("J61").Value = 7802.1665
$ws.Range("K61").Value = 5915
$ws.Range("L61").Value = 7802.1665
$ws.Range("M61").Value = -5703
$ws.Range("N61").Value = -8226.166499999999
$ws.Range("H74").Value = 19232112
$ws.Range("I74").Value = 22728546
$ws.Range("J74").Value = 1725
$ws.Range("K74").Value = 22728546
$ws.Range("L74").Value = 1725
$ws.Range("M74").Value = -22727672
$ws.Range("N74").Value = -3473
$ws.Range("H77").Value = 19232112
$ws.Range("I77").Value = 22728546
$ws.Range("J77").Value = 1725
$ws.Range("K77").Value = 113642730
$ws.Range("L77").Value = 8625
$ws.Range("M77").Value = -113638362
$ws.Range("N77").Value = -17361
$ws.Range("H97").Value = 2389.9092
$ws.Range("J97").Value = 3618.8
$ws.Range("L97").Value = 3618.8
$ws.Range("N97").Value = -4610.8
$ws.Range("H116").Value = 12675.791
$ws.Range("I116").Value = 16651.25
$ws.Range("K116").Value = 16651.25
$ws.Range("M116").Value = -14357.25
$ws.Range("H122").Value = 2717.2432
$ws.Range("I122").Value = 1591.2667
$ws.Range("K122").Value = 4773.800099999999
$ws.Range("M122").Value = -2323.800099999999
$ws.Range("H132").Value = 6857.7095
$ws.Range("I132").Value = 3243.6875
$ws.Range("K132").Value = 9731.0625
$ws.Range("M132").Value = -7201.0625
$ws.Range("H136").Value = 6238.514
$ws.Range("I136").Value = 5915
$ws.Range("J136").Value = 7802.1665
$ws.Range("K136").Value = 17745
$ws.Range("L136").Value = 23406.4995
$ws.Range("M136").Value = -15195
$ws.Range("N136").Value = -28506.4995
$ws.Range("H138").Value = 49999.5
$ws.Range("J138").Value = 49999.5
$ws.Range("L138").Value = 49999.5
$ws.Range("N138").Value = -60279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12675.791
$ws.Range("I3").Value = 16651.25
$ws.Range("K3").Value = 16651.25
$ws.Range("M3").Value = -16537.25
$ws.Range("H4").Value = 1300
$ws.Range("I4").Value = 1287.5
$ws.Range("K4").Value = 1287.5
$ws.Range("M4").Value = -1172.5
$ws.Range("H20").Value = 3596.3438
$ws.Range("I20").Value = 3820.3
$ws.Range("J20").Value = 3494.5454
$ws.Range("K20").Value = 3820.3
$ws.Range("L20").Value = 3494.5454
$ws.Range("M20").Value = -3573.3
$ws.Range("N20").Value = -3988.5454
$ws.Range("H94").Value = 1837
$ws.Range("I94").Value = 1744.125
$ws.Range("K94").Value = 1744.125
$ws.Range("M94").Value = -1293.125
$ws.Range("H105").Value = 3602.9333
$ws.Range("I105").Value = 3680.1667
$ws.Range("J105").Value = 3294
$ws.Range("K105").Value = 3680.1667
$ws.Range("L105").Value = 3294
$ws.Range("M105").Value = -1933.1667
$ws.Range("N105").Value = -6788
$ws.Range("H134").Value = 2387.3953
$ws.Range("I134").Value = 2273.3076
$ws.Range("J134").Value = 3499.75
$ws.Range("K134").Value = 6819.9228
$ws.Range("L134").Value = 10499.25
$ws.Range("M134").Value = -4284.9228
$ws.Range("N134").Value = -15569.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2362.3333
$ws.Range("I22").Value = 2362.3333
$ws.Range("K22").Value = 2362.3333
$ws.Range("M22").Value = -2012.3333
$ws.Range("H132").Value = 43013040
$ws.Range("I132").Value = 51283816
$ws.Range("J132").Value = 4999.6
$ws.Range("K132").Value = 153851448
$ws.Range("L132").Value = 14998.8
$ws.Range("M132").Value = -153848918
$ws.Range("N132").Value = -20058.8
$ws.Range("H134").Value = 3331.889
$ws.Range("I134").Value = 3331.889
$ws.Range("K134").Value = 9995.667000000001
$ws.Range("M134").Value = -7460.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 308.96155
$ws.Range("I2").Value = 95.63636
$ws.Range("J2").Value = 465.4
$ws.Range("K2").Value = 95.63636
$ws.Range("L2").Value = 465.4
$ws.Range("M2").Value = 17.36364
$ws.Range("N2").Value = -691.4
$ws.Range("H92").Value = 10750
$ws.Range("J92").Value = 10750
$ws.Range("L92").Value = 10750
$ws.Range("N92").Value = -14494
$ws.Range("H93").Value = 47834
$ws.Range("J93").Value = 47834
$ws.Range("L93").Value = 47834
$ws.Range("N93").Value = -51578
$ws.Range("H113").Value = 1479.4445
$ws.Range("I113").Value = 1199.4286
$ws.Range("K113").Value = 1199.4286
$ws.Range("M113").Value = 970.5714
$ws.Range("H122").Value = 338131.34
$ws.Range("I122").Value = 716675.1
$ws.Range("J122").Value = 6905.5
$ws.Range("K122").Value = 2150025.3
$ws.Range("L122").Value = 20716.5
$ws.Range("M122").Value = -2147575.3
$ws.Range("N122").Value = -25616.5
$ws.Range("H123").Value = 44764
$ws.Range("J123").Value = 44764
$ws.Range("L123").Value = 44764
$ws.Range("N123").Value = -49664
$ws.Range("H132").Value = 501887.25
$ws.Range("I132").Value = 572442.5600000001
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 1717327.68
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -1714797.68
$ws.Range("N132").Value = -29060
$ws.Range("H133").Value = 92492.5
$ws.Range("J133").Value = 92492.5
$ws.Range("L133").Value = 92492.5
$ws.Range("N133").Value = -102612.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 394.6
$ws.Range("I55").Value = 174.9375
$ws.Range("K55").Value = 174.9375
$ws.Range("M55").Value = -1.9375
$ws.Range("H68").Value = 5426.4
$ws.Range("I68").Value = 4722.846
$ws.Range("J68").Value = 9999.5
$ws.Range("K68").Value = 4722.846
$ws.Range("L68").Value = 9999.5
$ws.Range("M68").Value = -3973.846
$ws.Range("N68").Value = -11497.5
$ws.Range("H71").Value = 5426.4
$ws.Range("I71").Value = 4722.846
$ws.Range("J71").Value = 9999.5
$ws.Range("K71").Value = 23614.23
$ws.Range("L71").Value = 49997.5
$ws.Range("M71").Value = -19870.23
$ws.Range("N71").Value = -57485.5
$ws.Range("H132").Value = 5230.013
$ws.Range("I132").Value = 4376.4604
$ws.Range("K132").Value = 13129.3812
$ws.Range("M132").Value = -10599.3812
$ws.Range("H136").Value = 4364.769
$ws.Range("I136").Value = 2928.8333
$ws.Range("K136").Value = 8786.499899999999
$ws.Range("M136").Value = -6236.499899999999
$ws.Range("H139").Value = 111905
$ws.Range("J139").Value = 111905
$ws.Range("L139").Value = 111905
$ws.Range("N139").Value = -122185

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 29948
$ws.Range("I41").Value = 29946
$ws.Range("K41").Value = 29946
$ws.Range("M41").Value = -29556
$ws.Range("H113").Value = 514
$ws.Range("I113").Value = 349.4
$ws.Range("K113").Value = 1048.2
$ws.Range("M113").Value = 1121.8
$ws.Range("H126").Value = 83334210
$ws.Range("I126").Value = 142857920
$ws.Range("J126").Value = 1001
$ws.Range("K126").Value = 428573760
$ws.Range("L126").Value = 3003
$ws.Range("M126").Value = -428571290
$ws.Range("N126").Value = -7943
